$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DTVChannel")

# Append the six new "CUTVDisabledChannel_#" rows (18-23) under the existing
# CUTV-related rows, mirroring the style already used by the rows above
# (A12:B17, a matching 6-row block) by copying their formatting down onto
# the new rows.
$ws.Range("A12:B17").Copy()
$ws.Range("A18:B23").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$names  = @("CUTVDisabledChannel_1","CUTVDisabledChannel_2","CUTVDisabledChannel_3","CUTVDisabledChannel_4","CUTVDisabledChannel_5","CUTVDisabledChannel_6")
$values = @(1,2,4,6,7,9)

for ($i = 0; $i -lt 6; $i++) {
    $row = 18 + $i
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Make DTVChannel the active/selected sheet (was MiniEPGScreen before),
# with the selection left on A28 as in the authored workbook.
$ws.Activate()
$ws.Range("A28").Select()
